$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.390.85"
$ws.Range("E2").Value = "  +0.48%  "

$ws.Range("D3").Value = "1.890.08"
$ws.Range("E3").Value = "  -0.86%  "

$ws.Range("E4").Value = "  -0.72%  "

$ws.Range("E5").Value = "  -2.69%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.692"
$ws.Range("E6").Value = "  -0.39%  "

$ws.Range("E7").Value = "  -0.84%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "42.86"
$ws.Range("E8").Value = "  +2.85%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.354"
$ws.Range("E9").Value = "  -3.08%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "52.88"
$ws.Range("E10").Value = "  -0.13%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0739"
$ws.Range("E11").Value = "  -2.27%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0968"
$ws.Range("E12").Value = "  -1.65%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "13.04"
$ws.Range("E13").Value = "  -0.22%  "

$ws.Range("D14").Value = "2.166.64"
$ws.Range("E14").Value = "  -0.77%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.765"
$ws.Range("E15").Value = "  +2.35%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.95"
$ws.Range("E16").Value = "  -0.67%  "

$ws.Range("D17").Value = "1.873.39"
$ws.Range("E17").Value = "  -0.99%  "

$ws.Range("D18").Value = "35.519.44"
$ws.Range("E18").Value = "  +0.93%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "73.30"
$ws.Range("E19").Value = "  -1.07%  "

$ws.Range("D20").Value = "0.0₃0825"
$ws.Range("E20").Value = "  -1.54%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "245.23"
$ws.Range("E21").Value = "  +0.76%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "12.82"
$ws.Range("E22").Value = "  -2.19%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.97"
$ws.Range("E23").Value = "  -2.05%  "

$ws.Range("B24").Value = "Toncoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.65"
$ws.Range("E24").Value = "  +7.88%  "

$ws.Range("B25").Value = "Dai"
$ws.Range("C25").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("E25").Value = "  -0.74%  "

$ws.Range("E26").Value = "  -5.99%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "166.37"
$ws.Range("E27").Value = "  -0.20%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.51"
$ws.Range("E28").Value = "  -1.45%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.36"
$ws.Range("E29").Value = "  -0.98%  "

$ws.Range("E30").Value = "  -2.31%  "

$ws.Range("D31").Value = "4.128.46"
$ws.Range("E31").Value = "  -0.01%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.76"
$ws.Range("E32").Value = "  +9.21%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.26"
$ws.Range("E33").Value = "  -1.65%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0586"
$ws.Range("E34").Value = "  +0.46%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.19"
$ws.Range("E35").Value = "  -0.45%  "

$ws.Range("E36").Value = "  -0.82%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.82"
$ws.Range("E37").Value = "  -11.82%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.843"
$ws.Range("E38").Value = "  -2.57%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.99"
$ws.Range("E39").Value = "  -1.30%  "

$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0694"
$ws.Range("E40").Value = "  +6.32%  "

$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0219"
$ws.Range("E41").Value = "  +1.36%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "17.16"
$ws.Range("E42").Value = "  -1.93%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "97.62"
$ws.Range("E43").Value = "  -1.90%  "

$ws.Range("E44").Value = "  -3.01%  "

$ws.Range("D45").Value = "1.291.39"
$ws.Range("E45").Value = "  -3.58%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.32"
$ws.Range("E46").Value = "  -6.17%  "

$ws.Range("B47").Value = "Gas"
$ws.Range("C47").Value = "https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.59"
$ws.Range("E47").Value = "  +5.90%  "

$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0794"
$ws.Range("E48").Value = "  +5.81%  "

$ws.Range("B49").Value = "HuobiToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.41"
$ws.Range("E49").Value = "  -0.97%  "

$ws.Range("E50").Value = "  -0.53%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.27"
$ws.Range("E51").Value = "  -5.47%  "
